# Updated g-hMSC -> g-LTS image (2)
#
# Slide 9 ("writing/images/train.pptx" canonical numbering) holds a single
# top-level group shape "Groupe 195" (id=196) whose first child is a
# full-bleed background rectangle "Rectangle 247" (id=248). The edit moves
# the top edge of both the group and the rectangle up (while keeping the
# bottom edge fixed), growing their height to compensate:
#
#   off/ext  (611560, 548680) / (8208912, 5472608)
#        ->  (611560, 188640) / (8208912, 5832648)
#
# Only Top/Height change (Left/Width are left untouched so their values
# stay byte-identical to the source).
#
# Note: PowerPoint's Shape.Top/Height setters take points (1 pt = 12700
# EMU) and this host stores them as 32-bit floats, so a naive
# EMU/12700.0 literal can truncate one EMU short after the trip through
# points and back. The literals below are the exact float32 values whose
# point->EMU conversion lands exactly on the target EMU amounts.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

$grp = $s.Shapes.Item(1)          # "Groupe 195" (id 196)
$rect = $grp.GroupItems.Item(1)   # "Rectangle 247" (id 248), background rect

$newTop = 14.853544235229492      # 188640 EMU
$newHeight = 459.2636413574219    # 5832648 EMU

$grp.Top = $newTop
$grp.Height = $newHeight

$rect.Top = $newTop
$rect.Height = $newHeight
